$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D (shifts old D:K -> G:N)
$ws.Range("D1:F1").EntireColumn.Insert()

# Copy number formats/styles from the (now-shifted) old D:K block, i.e. G:N,
# onto the freshly inserted D:F columns so the new cells match the look of the
# columns they were cloned from (date format for header rows, number format for data).
$ws.Range("G5:N102").Copy()
$ws.Range("D5:F102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D:F columns with the newest quarter data
$newData = @{
    7 = @(43524, 43433, 43342)
    8 = @(5835000, 7913000, 8440000)
    9 = @(2971000, 3298000, 3289000)
    10 = @(2864000, 4615000, 5151000)
    12 = @(601000, 611000, 567000)
    13 = @(0, 0, 0)
    14 = @(134000, 19000, 27000)
    15 = @(0, 0, 0)
    17 = @(3961000, 4140000, 4062000)
    18 = @(1874000, 3773000, 4378000)
    20 = @(57000, 33000, 18000)
    21 = @(3244000, 5141000, 5681000)
    22 = @(27000, 33000, 50000)
    23 = @(1904000, 3773000, 4346000)
    24 = @(294000, 430000, 103000)
    25 = @(0, 0, 0)
    26 = @(1610000, 3343000, 4243000)
    27 = @(1605000, 3340000, 4242000)
    28 = @(0, 0, 0)
    29 = @(14000, -47000, 83000)
    30 = @(0, 0, 0)
    31 = @(0, 0, 0)
    32 = @(-57000, -33000, -18000)
    33 = @(1619000, 3293000, 4325000)
    34 = @(0, 0, 0)
    35 = @(1619000, 3293000, 4325000)
    38 = @(43524, 43433, 43342)
    41 = @(6353000, 4447000, 6506000)
    42 = @(1180000, 1116000, 296000)
    43 = @(4403000, 5418000, 5478000)
    44 = @(4390000, 3876000, 3595000)
    45 = @(224000, 182000, 164000)
    46 = @(16550000, 15039000, 16039000)
    47 = @(1614000, 1565000, 473000)
    48 = @(26204000, 24807000, 23672000)
    49 = @(1578000, 1584000, 1559000)
    50 = @(0, 0, 0)
    51 = @(0, 0, 0)
    52 = @(1541000, 1600000, 1633000)
    53 = @(0, 0, 0)
    54 = @(47487000, 44595000, 43376000)
    57 = @(1523000, 1683000, 1692000)
    58 = @(2634000, 398000, 859000)
    59 = @(3204000, 3108000, 3203000)
    60 = @(7361000, 5189000, 5754000)
    61 = @(3606000, 3736000, 3780000)
    62 = @(993000, 834000, 581000)
    63 = @(0, 0, 0)
    64 = @(0, 0, 0)
    65 = @(0, 0, 0)
    66 = @(12920000, 10726000, 11082000)
    68 = @(0, 0, 0)
    69 = @(0, 0, 0)
    70 = @(0, 0, 0)
    71 = @(0, 0, 0)
    72 = @(29364000, 27769000, 24395000)
    73 = @(0, 0, 0)
    74 = @(0, 0, 0)
    75 = @(0, 0, 0)
    76 = @(34567000, 33869000, 32294000)
    77 = @(0, 0, 0)
    80 = @(43524, 43433, 43342)
    81 = @(1619000, 3293000, 4325000)
    83 = @(1313000, 1335000, 1285000)
    84 = @(0, 0, 0)
    85 = @(0, 0, 0)
    86 = @(0, 0, 0)
    87 = @(0, 0, 0)
    88 = @(0, 0, 0)
    89 = @(3435000, 4810000, 5155000)
    91 = @(-2649000, -2700000, -2251000)
    92 = @(0, 0, 0)
    93 = @(0, 0, 0)
    94 = @(-2492000, -4427000, -2129000)
    96 = @(0, 0, 0)
    97 = @(0, 0, 0)
    98 = @(0, 0, 0)
    99 = @(0, 0, 0)
    100 = @(952000, -2435000, -3333000)
    101 = @(9000, -10000, -33000)
    102 = @(1904000, -2062000, -340000)
}
foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
}

# Data correction: row 91 (Net Borrowings) column J (was G pre-insert) changes
$ws.Cells.Item(91, 10).Value = -1265000

